$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Apply cell formatting (styles) by copying from representative
# "template" cells elsewhere in the sheet that already carry the desired
# style, then pasting formats-only onto the new Section 6 cells.

$ws.Range("A2").Copy()
$targets = @("A138", "A139", "A140", "A141", "A142", "A145", "A146", "A150", "A151", "A152", "A153", "A157", "A158", "A159", "A160", "A161", "A164", "A165", "A169", "A170", "A171", "A172", "A176", "A177")
foreach ($addr in $targets) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("B2").Copy()
$targets = @("B138", "B140", "B141", "B142", "B143", "B144", "B145", "B146", "B147", "B148", "B149", "B150", "B151", "C151", "B152", "B153", "B154", "B155", "B156", "B157", "B158", "C158", "B159", "B160", "B161", "B162", "B163", "B164", "B165", "B166", "B167", "B168", "B169", "B170", "C170", "B171", "B172", "B173", "B174", "B175", "B176", "B177", "C177")
foreach ($addr in $targets) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("A9").Copy()
$targets = @("A143", "A144", "A147", "A148", "A149", "A154", "A155", "A156", "A162", "A163", "A166", "A167", "A168", "A173", "A174", "A175")
foreach ($addr in $targets) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("B45").Copy()
$targets = @("B139")
foreach ($addr in $targets) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Step 2: Populate the text content (values) for the new Section 6 block.
$values = @{
    "A138" = "h1"
    "B138" = "頁面名稱 -6"
    "A139" = "p1"
    "B139" = "當初，蓬勃森林生長極短滿月之夜，它能夠讓人體驗前所未有的奇幻之旅。"
    "A140" = "h2"
    "B140" = "功能區塊 6.1"
    "A141" = "p"
    "B141" = "此外，「大師」滿足不同需求，適合使用。"
    "A142" = "ul"
    "B142" = "強大的代碼自動完成功能，提高開發效率。"
    "B143" = "內置語法檢查，幫助減少錯誤。"
    "B144" = "錯誤修正工具，輕鬆解決問題。"
    "A145" = "h3"
    "B145" = "單一功能名稱 6.1.1"
    "A146" = "ul"
    "B146" = "強大的代碼自動完成功能，提高開發效率。"
    "B147" = "內置語法檢查，幫助減少錯誤。"
    "B148" = "錯誤修正工具，輕鬆解決程式碼中的問題。"
    "B149" = "支持多種主流需求。"
    "B150" = "提供切換功能，方便使用。"
    "A151" = "img"
    "B151" = "img-1.png"
    "C151" = "容易上手。客製化設置，符合個性化需求。強大的幫助減少錯誤直觀的操作界面。"
    "A152" = "h3"
    "B152" = "單一功能名稱 6.1.2"
    "A153" = "ul"
    "B153" = "強大的代碼自動完成功能，提高開發效率。"
    "B154" = "內置語法檢查，幫助減少錯誤。"
    "B155" = "錯誤修正工具，輕鬆解決程式碼中的問題。"
    "B156" = "支持多種主流需求。"
    "B157" = "提供切換功能，方便使用。"
    "A158" = "img"
    "B158" = "img-1.png"
    "C158" = "容易上手。客製化設置，符合個性化需求。強大的幫助減少錯誤直觀的操作界面。"
    "A159" = "h2"
    "B159" = "功能區塊 6.2"
    "A160" = "p"
    "B160" = "此外，「大師」滿足不同需求，適合使用。"
    "A161" = "ul"
    "B161" = "強大的代碼自動完成功能，提高開發效率。"
    "B162" = "內置語法檢查，幫助減少錯誤。"
    "B163" = "錯誤修正工具，輕鬆解決問題。"
    "A164" = "h3"
    "B164" = "單一功能名稱 6.2.1"
    "A165" = "ul"
    "B165" = "強大的代碼自動完成功能，提高開發效率。"
    "B166" = "內置語法檢查，幫助減少錯誤。"
    "B167" = "錯誤修正工具，輕鬆解決程式碼中的問題。"
    "B168" = "支持多種主流需求。"
    "B169" = "提供切換功能，方便使用。"
    "A170" = "img"
    "B170" = "img-1.png"
    "C170" = "容易上手。客製化設置，符合個性化需求。強大的幫助減少錯誤直觀的操作界面。"
    "A171" = "h3"
    "B171" = "單一功能名稱 6.2.2"
    "A172" = "ul"
    "B172" = "強大的代碼自動完成功能，提高開發效率。"
    "B173" = "內置語法檢查，幫助減少錯誤。"
    "B174" = "錯誤修正工具，輕鬆解決程式碼中的問題。"
    "B175" = "支持多種主流需求。"
    "B176" = "提供切換功能，方便使用。"
    "A177" = "img"
    "B177" = "img-1.png"
    "C177" = "容易上手。客製化設置，符合個性化需求。強大的幫助減少錯誤直觀的操作界面。"
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
